$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest Adafruit IO reading as row 67.
$row = 67

$values = @(
    "2024-09-25T18:06:40Z",
    "temperature",
    "25",
    "N/A",
    "N/A",
    "N/A"
)

for ($col = 1; $col -le 6; $col++) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $values[$col - 1]
    $cell.Style = "Normal"
}
